$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.982.93'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.908.36'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7804'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.68'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3154'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.20'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06870'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07976'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = '1.904.45'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7398'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.194'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.69'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = '29.991.33'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.868'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -5.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.47'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007723'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9999'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').Value = '2.148.68'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.846'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '168.52'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.241'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1380'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +8.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.83'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.027'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.366'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.517'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.301'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05563'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.94%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.070'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.250'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7323'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.716'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01924'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.786'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.125'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4404'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.92'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8418'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.868'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.40'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.505'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.721'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '987.39'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +8.56%  '
$ws.Range('D50').Value = '2.055.84'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.19'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.18%  '
